# Generate Report for Handoff
# Updates the Overview, zh-cn and de-de sheets with refreshed handoff/handback
# status for the two files (7b09e590-... and e6081861-...).
#
# 7b09e590-...  -> now "Ready for handoff" (new handoff just went out)
# e6081861-...  -> "Handed back: in sync with en-US" with refreshed timestamp

$wb = $excel.ActiveWorkbook

$sevenB = "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
$e608    = "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"

$handedBack   = "Handed back: in sync with en-US"
$readyForHandoff = "Ready for handoff"

$handbackDate       = "2016-03-25 09:02:57"
$handoffDatetime     = "2016-03-25 09:02:48"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 -> e6081861 record, still handed back, refreshed datetime
$ov.Range("A2").Value = $e608
$ov.Range("B2").Value = $handedBack
$ov.Range("C2").Value = $handedBack
$ov.Range("D2").Value = $handbackDate

# Row 3 -> 7b09e590 record, now ready for handoff
$ov.Range("A3").Value = $sevenB
$ov.Range("B3").Value = $readyForHandoff
$ov.Range("C3").Value = $readyForHandoff
$ov.Range("D3").Value = $handbackDate

$ov.Hyperlinks.Item(1).TextToDisplay = $e608
$ov.Hyperlinks.Item(2).TextToDisplay = $sevenB

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhE608Handoff = "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf"
$zh7bHandoff   = "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf"

# Row 2 -> e6081861 record
$zh.Range("A2").Value = $e608
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $handedBack
$zh.Range("D2").Value = $zhE608Handoff
$zh.Range("E2").Value = $handoffDatetime
$zh.Range("F2").Value = $e608
$zh.Range("G2").Value = $zhE608Handoff
$zh.Range("H2").Value = "2016-03-25 09:01:39"
$zh.Range("J2").Value = "Include"

# Row 3 -> 7b09e590 record
$zh.Range("A3").Value = $sevenB
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $readyForHandoff
$zh.Range("D3").Value = $zh7bHandoff
$zh.Range("E3").Value = $handoffDatetime
$zh.Range("F3").Value = $sevenB
$zh.Range("G3").Value = $zh7bHandoff
$zh.Range("H3").Value = "2016-03-25 09:01:39"
$zh.Range("J3").Value = "Include"

$zh.Hyperlinks.Item(1).TextToDisplay = $e608
$zh.Hyperlinks.Item(2).TextToDisplay = $zhE608Handoff
$zh.Hyperlinks.Item(3).TextToDisplay = $e608
$zh.Hyperlinks.Item(4).TextToDisplay = $zhE608Handoff
$zh.Hyperlinks.Item(5).TextToDisplay = $sevenB
$zh.Hyperlinks.Item(6).TextToDisplay = $zh7bHandoff
$zh.Hyperlinks.Item(7).TextToDisplay = $sevenB
$zh.Hyperlinks.Item(8).TextToDisplay = $zh7bHandoff

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deE608Handoff = "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf"
$de7bHandoff   = "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf"

# Row 2 -> e6081861 record
$de.Range("A2").Value = $e608
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $handedBack
$de.Range("D2").Value = $deE608Handoff
$de.Range("E2").Value = $handbackDate
$de.Range("F2").Value = $e608
$de.Range("G2").Value = $deE608Handoff
$de.Range("H2").Value = "2016-03-25 09:01:59"
$de.Range("J2").Value = "Include"

# Row 3 -> 7b09e590 record
$de.Range("A3").Value = $sevenB
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $readyForHandoff
$de.Range("D3").Value = $de7bHandoff
$de.Range("E3").Value = $handbackDate
$de.Range("F3").Value = $sevenB
$de.Range("G3").Value = $de7bHandoff
$de.Range("H3").Value = "2016-03-25 09:01:59"
$de.Range("J3").Value = "Include"

$de.Hyperlinks.Item(1).TextToDisplay = $e608
$de.Hyperlinks.Item(2).TextToDisplay = $deE608Handoff
$de.Hyperlinks.Item(3).TextToDisplay = $e608
$de.Hyperlinks.Item(4).TextToDisplay = $deE608Handoff
$de.Hyperlinks.Item(5).TextToDisplay = $sevenB
$de.Hyperlinks.Item(6).TextToDisplay = $de7bHandoff
$de.Hyperlinks.Item(7).TextToDisplay = $sevenB
$de.Hyperlinks.Item(8).TextToDisplay = $de7bHandoff
